# update default image extension
# Insert a new "batch size" column (position 3) into the Run-tracking table,
# add 3 new run rows, and update the "transform"/"split" values accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- 1. Extend the style used by the existing table body (col L, s="1") onto the
#        new column M across every table row so the inserted column keeps the same
#        formatting as the rest of the table once the sheet grid reshuffles. ---
$ws.Range("L2:L30").Copy() | Out-Null
$ws.Range("M2:M30").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Write the new header row (A1:M1) ---
$headers = @("Run","model name","batch size","split","train folds","val folds","test folds","num classes","drop rate","in channels","Dataset","transform","image dir")
$hdrArr = New-Object 'object[,]' 1,13
for ($c = 0; $c -lt 13; $c++) { $hdrArr[0,$c] = $headers[$c] }
$ws.Range("A1:M1").Value2 = $hdrArr

# --- 3. Write the data rows (A2:M5) ---
# Row 2: existing resnet50_001 run, now carrying an explicit batch size of 32.
# Row 3: same run, re-exported with the new 320x224 transform.
# Row 4: same run/transform, but against the stricter balanced split.
# Row 5: new mobilenetv2_120d_001 run, strict split + 320x224 transform.
$dataRows = @(
    @("resnet50_001","resnet50",32,"full_balanced","0,1,2,3,4,5,6,7",8,9,5,0,3,"FrontViewDataset","create_transform((3,224,224))","C:\Users\Daniel\Documents\Data\Batch1"),
    @("resnet50_001","resnet50",32,"full_balanced","0,1,2,3,4,5,6,7",8,9,5,0,3,"FrontViewDataset","create_transform((3,320,224))","C:\Users\Daniel\Documents\Data\Batch1"),
    @("resnet50_001","resnet50",32,"strict_full_balanced","0,1,2,3,4,5,6,7",8,9,5,0,3,"FrontViewDataset","create_transform((3,320,224))","C:\Users\Daniel\Documents\Data\Batch1"),
    @("mobilenetv2_120d_001","mobilenetv2_120d",32,"strict_full_balanced","0,1,2,3,4,5,6,7",8,9,5,0,3,"FrontViewDataset","create_transform((3,320,224))","C:\Users\Daniel\Documents\Data\Batch1")
)
$dataArr = New-Object 'object[,]' 4,13
for ($r = 0; $r -lt 4; $r++) {
    for ($c = 0; $c -lt 13; $c++) {
        $dataArr[$r,$c] = $dataRows[$r][$c]
    }
}
$ws.Range("A2:M5").Value2 = $dataArr

# --- 4. Grow the table to cover the new column/rows; column names resync from the
#        header row automatically. ---
$lo.Resize($ws.Range("A1:M30"))

# --- 5. Column widths: best-fit widths for the now-shifted/renamed columns. ---
$ws.Columns.Item(3).ColumnWidth = 12.28515625
$ws.Columns.Item(4).ColumnWidth = 18.5703125
$ws.Columns.Item(5).ColumnWidth = 18.28515625
$ws.Columns.Item(9).ColumnWidth = 11.42578125
$ws.Columns.Item(10).ColumnWidth = 13.5703125
$ws.Columns.Item(11).ColumnWidth = 16.5703125
$ws.Columns.Item(12).ColumnWidth = 27.5703125
$ws.Columns.Item(13).ColumnWidth = 43.85546875

# --- 6. Selection cursor matches the post-edit save state. ---
$ws.Range("D19").Select()

Write-Host "Done"
